$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aug18")

# New round header (August 12, 2018 -> serial date 43324)
$ws.Range("A22").Value = 43324
$ws.Range("B22").Value = "Score"
$ws.Range("C22").Value = "Fairway"
$ws.Range("D22").Value = "GIR"
$ws.Range("E22").Value = "Putts"
$ws.Range("F22").Value = "Comment"

# Hole-by-hole data for the August 12 round
$ws.Range("A23").Value = "Hole 1"
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = "L"
$ws.Range("E23").Value = 0

$ws.Range("A24").Value = "Hole 2"
$ws.Range("B24").Value = 4
$ws.Range("E24").Value = 2

$ws.Range("A25").Value = "Hole 3"
$ws.Range("B25").Value = 4
$ws.Range("C25").Value = "R"
$ws.Range("E25").Value = 1

$ws.Range("A26").Value = "Hole 4"
$ws.Range("B26").Value = 7
$ws.Range("C26").Value = "L"
$ws.Range("E26").Value = 1

$ws.Range("A27").Value = "Hole 5"
$ws.Range("B27").Value = 3
$ws.Range("E27").Value = 1

$ws.Range("A28").Value = "Hole 6"
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = "S"
$ws.Range("E28").Value = 2

$ws.Range("A29").Value = "Hole 7"
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = "S"
$ws.Range("E29").Value = 2

$ws.Range("A30").Value = "Hole 8"
$ws.Range("B30").Value = 4
$ws.Range("E30").Value = 2

$ws.Range("A31").Value = "Hole 9"
$ws.Range("B31").Value = 5
$ws.Range("C31").Value = "S"
$ws.Range("E31").Value = 2

$ws.Range("A32").Value = "Hole 10"
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = "S"
$ws.Range("E32").Value = 1

$ws.Range("A33").Value = "Hole 11"
$ws.Range("B33").Value = 4
$ws.Range("E33").Value = 2

$ws.Range("A34").Value = "Hole 12"
$ws.Range("B34").Value = 4
$ws.Range("C34").Value = "S"
$ws.Range("E34").Value = 1

$ws.Range("A35").Value = "Hole 13"
$ws.Range("B35").Value = 6
$ws.Range("C35").Value = "R"
$ws.Range("E35").Value = 4

$ws.Range("A36").Value = "Hole 14"
$ws.Range("B36").Value = 5
$ws.Range("C36").Value = "S"
$ws.Range("E36").Value = 2

$ws.Range("A37").Value = "Hole 15"
$ws.Range("B37").Value = 3
$ws.Range("E37").Value = 1

$ws.Range("A38").Value = "Hole 16"
$ws.Range("B38").Value = 5
$ws.Range("C38").Value = "R"
$ws.Range("E38").Value = 3

$ws.Range("A39").Value = "Hole 17"
$ws.Range("B39").Value = 5
$ws.Range("C39").Value = "L"
$ws.Range("E39").Value = 3

$ws.Range("A40").Value = "Hole 18"
$ws.Range("B40").Value = 4
$ws.Range("C40").Value = "S"
$ws.Range("E40").Value = 2

# Totals row
$ws.Range("B41").Formula = "=SUM(B23:B40)"
$ws.Range("E41").Formula = "=SUM(E23:E40)"
